# "correcion base de datos" -- fix labels/values in the mortality table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "Hemorragia en el embarazo..." -> shortened to "Hemorragia"
$ws.Range("B6").Value = "Hemorragia"

# Row 7: "Complicaciones predominantes, relacionadas..." -> shortened to "Complicaciones predominantes"
$ws.Range("B7").Value = "Complicaciones predominantes"

# Row 9: label swapped to "Complicaciones relacionadas con el puerperio" and NA values replaced with 0
$ws.Range("B9").Value = "Complicaciones relacionadas con el puerperio"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# Row 11: long label shortened, dropping the leading "Muerte materna debida a"
$ws.Range("B11").Value = "Cualquier causa obstétrica que ocurre después de 42 días pero antes de un año del parto"

# Column B widened to fit the (now shorter, but still present) long labels
# (Excel's stored <col width> = ColumnWidth + 5/6, so back the offset out here
#  to land exactly on width="41" / width="12" in the saved XML.)
$ws.Range("B:B").ColumnWidth = 40.16666666666667
$ws.Range("C:D").ColumnWidth = 11.16666666666667

# Enable iterative calculation (calcPr iterate="1")
$excel.Iteration = $true

# Selection moved to G11 in the saved file
$ws.Range("G11").Select() | Out-Null
